$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 133

for ($r = 3; $r -le $lastRow; $r++) {
    if ($r -eq 11 -or $r -eq 52 -or $r -eq 93) {
        continue
    }
    $ws.Range("K$r").Value = 2
    $ws.Range("L$r").Value = 1
    $ws.Range("M$r").Value = 6
    $ws.Range("N$r").Value = 3
    $ws.Range("O$r").Value = 3
    $ws.Range("P$r").Value = 6
    $ws.Range("Q$r").Value = 9
}

$ws.Range("Q11").Select()
